# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) on Hoja1 lists YYMM period codes for rows
# 16-59. The database refresh re-sorts these periods from oldest-first to
# newest-first (ascending instead of descending), i.e. the list is reversed:
#   before: 2003, 2002, 2001, 1912, ... , 1609, 1608
#   after : 1608, 1609, 1610, 1611, ... , 2002, 2003

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @(
    "1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 5).Value = $periodos[$i]
}
